$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "release/8.0.8"
$ws.Range("B11").Value = "X"
$ws.Range("C11").Value = "X"
$ws.Range("D11").Value = "X"
$ws.Range("E11").Value = "X"
